# Regenerate the "K" column (G) for the save_data sheet.
# The sheet previously stored a "Strike#"-style raw value in column G;
# this recomputes/rewrites it as the new "K" statistic (s_vals) for
# every data row (rows 2-48) in one vectorized range write.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New K values, one per row for G2:G48 (in row order).
$s_vals = @(
    3, 5, 4, 1, 2, 0, 4, 2, 1, 0,
    0, 2, 1, 2, 3, 0, 0, 1, 3, 3,
    0, 0, 4, 2, 0, 3, 1, 2, 1, 0,
    0, 1, 0, 0, 2, 2, 1, 0, 3, 0,
    3, 1, 1, 1, 1, 3, 3
)

$firstRow = 2
$lastRow = $firstRow + $s_vals.Count - 1

# Build a vertical (N x 1) array to assign to the whole column range at once.
$arr = New-Object 'object[,]' $s_vals.Count, 1
for ($i = 0; $i -lt $s_vals.Count; $i++) {
    $arr[$i, 0] = $s_vals[$i]
}

$targetRange = $ws.Range("G$firstRow`:G$lastRow")
$targetRange.Value = $arr
